$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update effect_id (column C) values for existing rows ("일반 몬스터 spell_effect 변경") ---
$ws.Range("C4").Value = 10103
$ws.Range("C5").Value = 10104
$ws.Range("C6").Value = 10105
$ws.Range("C7").Value = 10101
$ws.Range("C8").Value = 20101
$ws.Range("C9").Value = 30201
$ws.Range("C10").Value = 10101
$ws.Range("C11").Value = 20101
$ws.Range("C12").Value = 30201
$ws.Range("C13").Value = 10104
$ws.Range("C14").Value = 30205
$ws.Range("C15").Value = 20102
$ws.Range("C16").Value = 10101
$ws.Range("C17").Value = 20101
$ws.Range("C18").Value = 30201
$ws.Range("C19").Value = 10101
$ws.Range("C20").Value = 20101
$ws.Range("C21").Value = 30201
$ws.Range("C22").Value = 10103
$ws.Range("C23").Value = 10205
$ws.Range("C24").Value = 30203
$ws.Range("C25").Value = 10205
$ws.Range("C26").Value = 20102
$ws.Range("C27").Value = 10101
$ws.Range("C28").Value = 20101
$ws.Range("C29").Value = 30201
$ws.Range("C30").Value = 10101
$ws.Range("C31").Value = 20101
$ws.Range("C32").Value = 30201
$ws.Range("C33").Value = 10101
$ws.Range("C34").Value = 30201
$ws.Range("C35").Value = 10101
$ws.Range("C36").Value = 20101
$ws.Range("C37").Value = 30201
$ws.Range("C38").Value = 10101
$ws.Range("C39").Value = 20101
$ws.Range("C40").Value = 30201

# --- Append new rows 41-46 for 몬스터_전사_엘리트_LV1 ("몬스터 spell_effect 변경") ---
# Seed formatting for the new rows by copying the style of the last existing data row (row 40)
$ws.Range("A40:G40").Copy()
$ws.Range("A41:G46").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A41").Value = 201010101
$ws.Range("B41").Value = "몬스터_전사_엘리트_LV1"
$ws.Range("C41").Value = 10106
$ws.Range("D41").Value = "대미지 타입 효과"
$ws.Range("E41").Value = 3
$ws.Range("F41").Value = "none"

$ws.Range("A42").Value = 201010101
$ws.Range("B42").Value = "몬스터_전사_엘리트_LV1"
$ws.Range("C42").Value = 10107
$ws.Range("D42").Value = "대미지 타입 효과"
$ws.Range("E42").Value = 4
$ws.Range("F42").Value = "none"

$ws.Range("A43").Value = 201010101
$ws.Range("B43").Value = "몬스터_전사_엘리트_LV1"
$ws.Range("C43").Value = 10108
$ws.Range("D43").Value = "대미지 타입 효과"
$ws.Range("E43").Value = 5
$ws.Range("F43").Value = "none"

$ws.Range("A44").Value = 201010201
$ws.Range("B44").Value = "몬스터_전사_엘리트_LV1"
$ws.Range("C44").Value = 10108
$ws.Range("D44").Value = "대미지 타입 효과"
$ws.Range("E44").Value = 3
$ws.Range("F44").Value = "none"

$ws.Range("A45").Value = 201010201
$ws.Range("B45").Value = "몬스터_전사_엘리트_LV1"
$ws.Range("C45").Value = 10109
$ws.Range("D45").Value = "대미지 타입 효과"
$ws.Range("E45").Value = 4
$ws.Range("F45").Value = "none"

$ws.Range("A46").Value = 201010201
$ws.Range("B46").Value = "몬스터_전사_엘리트_LV1"
$ws.Range("C46").Value = 10110
$ws.Range("D46").Value = "대미지 타입 효과"
$ws.Range("E46").Value = 5
$ws.Range("F46").Value = "none"

# --- Column B width / bestfit adjustment for the new, wider monster name ---
$ws.Columns.Item(2).ColumnWidth = 17.4

# --- Restore the active selection exactly as in the authored workbook ---
$ws.Range("H20").Select()